# Add a new "SydneyQa2019" record row to the "SignIn" sheet (xl/worksheets/
# sheet2.xml, which is Worksheets.Item(2)), mirroring the existing
# "SydneyQa2018" row: a hyperlinked e-mail in column B and the password
# text in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# New row of data (row 3): B3 = same e-mail as B2, C3 = new password string.
$ws.Range("B3").Value = "mvpstudio.qa@gmail.com"
$ws.Range("C3").Value = "SydneyQa2019"

# Give B3 a mailto hyperlink just like B2, then restore the shared
# "Hyperlink" cell style (Add() applies its own style variant otherwise).
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:mvpstudio.qa@gmail.com")
$ws.Range("B3").Style = "Hyperlink"

# Leave the sheet's selection where the author ended up after editing.
$ws.Range("C6").Select()
